# Two match records in the data table got their row content (everything
# except the running index in column A) swapped during the base update:
#   - row 129  <->  row 131
#   - row 224  <->  row 225
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(@(129, 131), @(224, 225))

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    # Columns B (2) through AD (30) hold the match record; column A is the
    # running index and stays put.
    $rng1 = $ws.Range($ws.Cells.Item($row1, 2), $ws.Cells.Item($row1, 30))
    $rng2 = $ws.Range($ws.Cells.Item($row2, 2), $ws.Cells.Item($row2, 30))

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
